$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New data entries added below the existing rows (rows 4, 5, 8, 12).
# Cell order below matches the order new shared-string entries were
# introduced so the rebuilt sharedStrings table lines up with the target.
$ws.Range("E4").Value = "Testyantra"

$ws.Range("A5").Value = "Yogendra"
$ws.Range("E5").Value = "Testyantra"

$ws.Range("B8").Value = "Mr."
$ws.Range("C8").Value = "raja"
$ws.Range("D8").Value = "nalla"

$ws.Range("C5").Value = "Biotechnology"

$ws.Range("B12").Value = "Employee"
$ws.Range("C12").Value = "Education"

# Match the final selection shown in the saved workbook
$ws.Range("C12").Select()
